$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.710899999999995
$ws.Range("B4").Value = 9.040799999999996
$ws.Range("D4").Value = -7.231599999999994
$ws.Range("C5").Value = -12.38390000000001
$ws.Range("D5").Value = -7.679099999999997
$ws.Range("A8").Value = -22.10340000000001
$ws.Range("C8").Value = -13.1068
$ws.Range("A10").Value = -22.1715
$ws.Range("B11").Value = 5.810199999999997
$ws.Range("A12").Value = -21.86859999999999
$ws.Range("B12").Value = 5.9041
$ws.Range("C12").Value = -11.3882
$ws.Range("C13").Value = -11.89509999999999
$ws.Range("B15").Value = 5.627399999999994
$ws.Range("C15").Value = -14.55429999999999
$ws.Range("B17").Value = 5.286099999999998
$ws.Range("A18").Value = -22.29410000000001
$ws.Range("D20").Value = -8.717699999999986
$ws.Range("C21").Value = -13.2741
$ws.Range("D23").Value = -8.329499999999996
$ws.Range("A25").Value = -21.69159999999999
$ws.Range("C25").Value = -12.7091
$ws.Range("B26").Value = 5.267399999999998
$ws.Range("D26").Value = -7.584700000000001
$ws.Range("B27").Value = 5.474099999999993
$ws.Range("B28").Value = 6.014299999999999
$ws.Range("B32").Value = 7.9447
$ws.Range("C32").Value = -12.2593
$ws.Range("D34").Value = -7.802900000000007
$ws.Range("C36").Value = -13.49500000000002
$ws.Range("A37").Value = -21.77740000000001
$ws.Range("B37").Value = 5.759399999999997
$ws.Range("C38").Value = -12.41379999999999
$ws.Range("D39").Value = -7.797899999999994
$ws.Range("D40").Value = -8.738199999999988
$ws.Range("B41").Value = 8.897899999999998
$ws.Range("C41").Value = -13.11320000000002
$ws.Range("D41").Value = -8.238999999999994
$ws.Range("D42").Value = -8.328899999999994
$ws.Range("B47").Value = 6.287399999999998
$ws.Range("D47").Value = -7.813700000000001
$ws.Range("C50").Value = -14.04729999999999
$ws.Range("B51").Value = 5.106300000000002
$ws.Range("C52").Value = -12.3174
$ws.Range("D52").Value = -7.519299999999998
$ws.Range("A55").Value = -22.3563
$ws.Range("C59").Value = -12.7715
$ws.Range("D60").Value = -8.643399999999994
$ws.Range("D62").Value = -8.757799999999989
$ws.Range("B65").Value = 5.804899999999999
$ws.Range("C67").Value = -11.48659999999999
$ws.Range("A68").Value = -21.481
$ws.Range("D70").Value = -7.008899999999994
$ws.Range("D72").Value = -7.280400000000004
$ws.Range("B73").Value = 9.2799
$ws.Range("A77").Value = -20.6139
$ws.Range("A78").Value = -19.83029999999998
$ws.Range("A79").Value = -20.64619999999999
$ws.Range("A80").Value = -19.3399
$ws.Range("A81").Value = -22.18890000000001
$ws.Range("A82").Value = -21.58310000000001
$ws.Range("D83").Value = -8.983400000000001
$ws.Range("A84").Value = -22.06800000000003
$ws.Range("B84").Value = 4.838499999999999
$ws.Range("C84").Value = -12.49219999999999
$ws.Range("B85").Value = 5.3119
$ws.Range("C86").Value = -13.63189999999999
$ws.Range("C88").Value = -13.159
$ws.Range("B89").Value = 4.897699999999997
$ws.Range("C89").Value = -14.2175
$ws.Range("B93").Value = 5.405999999999998
$ws.Range("B95").Value = 5.7205
$ws.Range("C95").Value = -12.98360000000001
$ws.Range("B98").Value = 7.824299999999998
$ws.Range("B99").Value = 5.868199999999998
$ws.Range("A101").Value = -21.74549999999999
$ws.Range("B101").Value = 5.691699999999996
$ws.Range("A102").Value = -22.05260000000002
$ws.Range("B102").Value = 5.200899999999998
$ws.Range("C105").Value = -13.3217
